# "camera can catch QR-code"
# The QR scan log sheet grew new rows: two more people (Елена Шалаева,
# МАКСИМ Вихров) plus a run of additional timestamps for Иван Кизикин.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A ("ФИО" / name) for every data row (2-15), existing rows first,
# then the newly added ones - mirrors the order the log rows were scanned.
$names = @(
    "Иван Кизикин",
    "Андрей Фокин",
    "Елена  Шалаева ",
    "МАКСИМ Вихров",
    "Иван Кизикин",
    "Иван Кизикин",
    "Иван Кизикин",
    "Иван Кизикин",
    "Иван Кизикин",
    "Иван Кизикин",
    "Иван Кизикин",
    "Иван Кизикин",
    "Иван Кизикин",
    "Иван Кизикин"
)

# Column B ("Время" / timestamp) for every data row (2-15), same order.
$times = @(
    "2023-06-27 21:23:11",
    "2023-07-17 12:19:35",
    "2023-07-20 11:13:08",
    "2023-07-23 20:25:07",
    "2023-07-24 16:31:13",
    "2023-07-24 16:33:50",
    "2023-07-24 16:36:15",
    "2023-07-24 16:39:18",
    "2023-07-24 16:43:53",
    "2023-07-24 16:48:59",
    "2023-07-24 16:56:30",
    "2023-07-24 17:00:39",
    "2023-07-24 17:03:29",
    "2023-07-24 17:05:15"
)

for ($i = 0; $i -lt $names.Count; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $names[$i]
}

for ($i = 0; $i -lt $times.Count; $i++) {
    $ws.Cells.Item($i + 2, 2).Value = $times[$i]
}
